$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.653.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.964.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.19'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0814'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.252.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.829'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.966.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.535.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0859'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.43%  '
$ws.Range('E26').Value = '  +7.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.17'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.119'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0618'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.25'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.00%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.97'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.22%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.95%  '
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0989'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +1.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.08'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.361.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.81'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.84'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.143.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.85'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.50%  '
